$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'302.63"
$ws.Range("D3").Value = "'35.65"
$ws.Range("E3").Value = "'-1.58%"
$ws.Range("E4").Value = "'-1.15%"
$ws.Range("D5").Value = "'0.07899"
$ws.Range("E5").Value = "'-2.79%"
$ws.Range("D6").Value = "'1.853"
$ws.Range("E6").Value = "'-4.66%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'4.106"
$ws.Range("E7").Value = "'-1.81%"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "'7.784"
$ws.Range("E8").Value = "'0.01%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9192"
$ws.Range("E9").Value = "'-1.27%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1344"
$ws.Range("E10").Value = "'-3.69%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1881"
$ws.Range("E11").Value = "'-2.17%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09013"
$ws.Range("E12").Value = "'-2.85%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03457"
$ws.Range("E13").Value = "'0.90%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09789"
$ws.Range("E14").Value = "'-0.75%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001407"
$ws.Range("E15").Value = "'-1.24%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006101"
$ws.Range("E16").Value = "'6.20%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.722"
$ws.Range("E17").Value = "'3.23%"
$ws.Range("D18").Value = "'3.299"
$ws.Range("E18").Value = "'10.99%"
$ws.Range("D19").Value = "'0.3438"
$ws.Range("D20").Value = "'5.173"
$ws.Range("E20").Value = "'5.45%"
$ws.Range("D23").Value = "'0.04412"
$ws.Range("E23").Value = "'-2.43%"
$ws.Range("E24").Value = "'1.64%"
$ws.Range("D25").Value = "'0.004602"
$ws.Range("E25").Value = "'-5.49%"
$ws.Range("E26").Value = "'4.86%"
$ws.Range("D27").Value = "'0.0004443"
$ws.Range("E27").Value = "'0.03%"
$ws.Range("E39").Value = "'-3.82%"
$ws.Range("D40").Value = "'0.05267"
$ws.Range("E40").Value = "'6.58%"
$ws.Range("D41").Value = "'0.007611"
$ws.Range("E41").Value = "'-0.62%"
$ws.Range("D42").Value = "'0.01016"
$ws.Range("E42").Value = "'-0.58%"
$ws.Range("E43").Value = "'-2.97%"
$ws.Range("D44").Value = "'0.002161"
$ws.Range("E44").Value = "'2.87%"
$ws.Range("D45").Value = "'0.01013"
$ws.Range("E45").Value = "'-11.85%"
$ws.Range("D46").Value = "'0.00006142"
$ws.Range("E46").Value = "'-4.89%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'0.02%"
$ws.Range("D48").Value = "'65.22"
$ws.Range("E48").Value = "'0.85%"
$ws.Range("D49").Value = "'0.001659"
$ws.Range("E49").Value = "'39.36%"
$ws.Range("E50").Value = "'0.02%"
$ws.Range("E51").Value = "'0.02%"
